$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"   = 7.93
    "B6"   = 5.906000000000001
    "B7"   = 6.025
    "B16"  = 5.389
    "B20"  = 6.891
    "B28"  = 6.309
    "B29"  = 5.615
    "B32"  = 6.845999999999999
    "B40"  = 9.344000000000001
    "B46"  = 6.627
    "B51"  = 5.765000000000001
    "B52"  = 5.825
    "B57"  = 5.189
    "B59"  = 4.539
    "B62"  = 5.261
    "B66"  = 5.115
    "B73"  = 7.066
    "B74"  = 9.204000000000001
    "B92"  = 5.882
    "B100" = 5.848999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
